$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 46, pushing existing rows 46:56 down to 47:57.
$ws.Rows.Item(46).Insert()

# Populate the newly inserted row 46 with the new weekly record.
$ws.Cells.Item(46, 1).Value = 4
$ws.Cells.Item(46, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(46, 3).Value = "Los Lagos"
$ws.Cells.Item(46, 4).Value = 44663
$ws.Cells.Item(46, 5).Value = 10
$ws.Cells.Item(46, 6).Value = 100112043
$ws.Cells.Item(46, 7).Value = "Pepino dulce"
$ws.Cells.Item(46, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(46, 9).Value = "Primera"
$ws.Cells.Item(46, 10).Value = 90
$ws.Cells.Item(46, 11).Value = 18000
$ws.Cells.Item(46, 12).Value = 18000
$ws.Cells.Item(46, 13).Value = 18000
$ws.Cells.Item(46, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(46, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(46, 16).Value = 1000
$ws.Cells.Item(46, 17).Value = 18
$ws.Cells.Item(46, 18).Value = "Hortaliza"
